# Commit: "limiting the dragging of images only to the visible area"
#
# Moves the "Align the images vertically centered..." paragraph to just
# after "User should be able to see the jumbled images." (followed by a
# blank paragraph before it, and two blank paragraphs after it, before
# "To Do"). Moves the _GoBack bookmark paragraph to just after "To Do"
# (followed by two blank paragraphs before "Prevent the images..."), and
# removes the stale <w:lastRenderedPageBreak/> marker from the
# "Settings will have ->" paragraph.

$d = $word.ActiveDocument

function Find-ParaByText($text) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $ok = $rng.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Find failed for: $text"
    }
    return $rng.Paragraphs(1)
}

$alignText = "Align the images vertically centered depending on the number of rows."

# --- Step 1: remove the two paragraphs whose *position* changes (the
#     Align text paragraph and the bookmark paragraph), plus the three
#     blank paragraphs that used to sit between "To Do" and
#     "Prevent the images..." -- all five are rebuilt fresh below, in
#     their new spots.

$pAlign = Find-ParaByText($alignText)
$pAlign.Range.Delete() | Out-Null

$pToDoHost = Find-ParaByText("To Do")
$pBookmark = $pToDoHost.Previous()
$pBookmark.Range.Delete() | Out-Null

$pToDo = Find-ParaByText("To Do")
$pToDo.Next().Range.Delete() | Out-Null
$pToDo = Find-ParaByText("To Do")
$pToDo.Next().Range.Delete() | Out-Null
$pToDo = Find-ParaByText("To Do")
$pToDo.Next().Range.Delete() | Out-Null

# --- Step 2: rebuild the block right after "jumbled images":
#     <blank> / Align-text / <blank> / <blank>
#     Built forward with InsertParagraphAfter so the new paragraphs
#     inherit the plain (unformatted) paragraph mark of "jumbled images"
#     rather than any heading formatting.

$pJumbled = Find-ParaByText("User should be able to see the jumbled images.")
$pJumbled.Range.InsertParagraphAfter() | Out-Null

$n1 = $pJumbled.Next()
$n1.Range.InsertParagraphAfter() | Out-Null

$n2 = $n1.Next()
$n2.Range.Text = $alignText

$n2.Range.InsertParagraphAfter() | Out-Null
$n3 = $n2.Next()
$n3.Range.InsertParagraphAfter() | Out-Null

# --- Step 3: rebuild the block right after "To Do":
#     _GoBack-bookmark-blank / <blank> / <blank>
#     Built backward with InsertParagraphBefore, anchored on
#     "Prevent the images..." (plain formatting) so the new paragraphs
#     don't inherit the "To Do" heading's bold/underline/size formatting.
#     The anchor paragraph is re-located via Find before each insert,
#     since InsertParagraphBefore repositions the paragraph it is called
#     on (earlier references to it go stale).

$pPrevent = Find-ParaByText("Prevent the images from being dragged out of the area of interest")
$pPrevent.Range.InsertParagraphBefore() | Out-Null

$pPrevent = Find-ParaByText("Prevent the images from being dragged out of the area of interest")
$pPrevent.Range.InsertParagraphBefore() | Out-Null

$pPrevent = Find-ParaByText("Prevent the images from being dragged out of the area of interest")
$pPrevent.Range.InsertParagraphBefore() | Out-Null

# The bookmark host is the (now blank) paragraph immediately after "To Do".
$pToDo = Find-ParaByText("To Do")
$bmHost = $pToDo.Next()

# Bookmarks.Add leaks bookmarkEnd into the next paragraph when given a
# range that exactly spans an empty paragraph (start==paragraph start,
# end==paragraph end/next-paragraph boundary). Work around it by
# temporarily typing a placeholder character, bookmarking exactly that
# character (a range that does not touch the paragraph-end boundary),
# then deleting the placeholder -- leaving the bookmark cleanly inside
# its own paragraph.
$bmHost.Range.Text = "X"
$bmSub = $bmHost.Range.Duplicate()
$bmSub.Find.Execute("X", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Bookmarks.Add("_GoBack", $bmSub) | Out-Null

$bmClear = $bmHost.Range.Duplicate()
$bmClear.Find.Execute("X", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

$bmHost.Range.InsertParagraphAfter() | Out-Null
$n4 = $bmHost.Next()
$n4.Range.InsertParagraphAfter() | Out-Null

# --- Step 4: drop the stale lastRenderedPageBreak marker on the
#     "Settings will have ->" paragraph by rewriting its run text (a
#     plain re-assignment of the same text rebuilds the run cleanly,
#     without the obsolete rendering-hint element, while leaving the
#     paragraph's own identity/attributes untouched).

$pSettings = Find-ParaByText("Settings will have ->")
$pSettings.Range.Text = "Settings will have ->"
